# Update ticket-count / price figures for the two affected sheets
# (展览 and 全部类型) per the commit's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 195
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F4").Value = 115
$ws1.Range("F5").Value = 645

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 195
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F5").Value = 115
$ws4.Range("F6").Value = 645
